$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store numbers/percentages as
# plain text rather than numeric cell values. Format each refreshed cell
# as Text before assigning its new value so it is written back the same
# way (e.g. "307.52"), instead of Excel auto-converting it to a
# number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "307.52"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-2.27%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "37.47"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-4.94%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.101"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-0.22%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07872"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-3.99%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.982"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.04%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "4.341"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.57%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "8.226"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.29%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-5.73%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9260"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "-0.61%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1277"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-9.71%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1894"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.94%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08892"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-2.23%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03425"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-3.35%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09752"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.58%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001395"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-0.60%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.005982"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.22%"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "1,779.46%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.581"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-2.19%"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.79%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1283"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-1.47%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.009"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.78%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.2495"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2.01%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04323"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.02%"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.03%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004598"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.01%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "176.76%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02291"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.39%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05011"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "-5.13%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007508"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.54%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009898"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "1.07%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1353"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-1.71%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.008020"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-18.12%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006521"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.51%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.36%"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.003003"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "8.71%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001202"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.32%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002104"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.36%"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002004"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.36%"
